# Generate Report for Handoff
# Updates the localization-status workbook: files that were previously
# "Handed back" are now "Ready for handoff" again, with refreshed handoff
# timestamps, and an error detail noting the handback file for
# e87b3ebd-18b1-49cf-b232-fe0371daea31 is stale.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/7ed0b58048952ab709cd959dea829061846689bb/e2e/e87b3ebd-18b1-49cf-b232-fe0371daea31.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/f51464ed22b4c674dc026a8620e177a323fe7a37/e2e/e87b3ebd-18b1-49cf-b232-fe0371daea31.md."

# --- Overview sheet: summary row for e87b3ebd-18b1-49cf-b232-fe0371daea31.md ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusReady
$overview.Range("F3").Value = $statusReady
$overview.Range("G3").Value = "2016-08-12 18:59:52"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")

# Widen the Error Detail column
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# Row 2: be89b7da-88b6-4d30-a0f4-ec864ee83e08.md
$zhcn.Range("C2").Value = $statusReady

# Row 3: e87b3ebd-18b1-49cf-b232-fe0371daea31.md
$zhcn.Range("C3").Value = $statusReady
$zhcn.Range("H3").Value = "2016-08-12 18:59:45"
$zhcn.Range("P3").Value = $errorDetail

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")

# Widen the Error Detail column
$dede.Columns.Item(16).ColumnWidth = 39.17

# Row 2: be89b7da-88b6-4d30-a0f4-ec864ee83e08.md
$dede.Range("C2").Value = $statusReady

# Row 3: e87b3ebd-18b1-49cf-b232-fe0371daea31.md
$dede.Range("C3").Value = $statusReady
$dede.Range("H3").Value = "2016-08-12 18:59:52"
$dede.Range("P3").Value = $errorDetail
